# Update marksheet: correct marks and total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row total-right value: 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row right value: 60 -> 100
$ws.Range("B12").Value = 100

# "Total" row Max column (Correct/Total marks text): 56/84 -> 100/140
$ws.Range("E12").Value = "100/140"

$wb.Save()
